# Add validations for denial_reasons and denial_reasons_ff fields (#14)
#
# The "action_taken_date" column (K) on the "invalid" sheet is removed, which
# shifts the app_recipient / amount_approved / amount_applied_for_flag /
# amount_applied_for columns one place to the left (L->K, M->L, N->M, O->N).
# Two brand-new columns are then appended for denial_reasons (O) and
# denial_reasons_ff (P), each with header + per-row sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Remove the action_taken_date column entirely - this shifts app_recipient,
# amount_approved, amount_applied_for_flag and amount_applied_for one
# column to the left (K<-L, L<-M, M<-N, N<-O).
$ws.Columns.Item(11).Delete()

# --- New column O: denial_reasons ---
$ws.Range("O1").Value = "denial_reasons"
$ws.Range("O2").Value = 1
$ws.Range("O3").Value = 101
$ws.Range("O4").Value = 1
$ws.Range("O5").Value = 102
$ws.Range("O6").Value = "1;2;999"
$ws.Range("O7").Value = 999
$ws.Range("O8").Value = "1;2;3;101;4"
$ws.Range("O9").Value = 1
$ws.Range("O10").Value = "1;3;4;999"
$ws.Range("O11").Value = 977

# --- New column P: denial_reasons_ff ---
$ws.Range("P1").Value = "denial_reasons_ff"
$ws.Range("P9").Value = "abc123"

# Leave the remaining P column cells (P2:P8, P10:P11) blank, matching the
# sparsely populated free-form text field in the source data.

# Update the sheet view to match the authored selection/scroll position.
$ws.Range("O5").Select()
$ws.Application.ActiveWindow.ScrollColumn = 6
